# trials.xlsx edit — "added trials thru 80"
#
# The source workbook tracks vehicle-accident trials. Every real trial row
# (e.g. row 108 = trial "54") is immediately followed by a companion
# "<n>trial" row (row 109 = "54trial") that holds the *re-reviewed* /
# corrected data for that same trial. Many of those companion rows were
# still blank; this commit fills in the ones for trials up through 80,
# plus a couple of small corrections to already-filled cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 109 ("54trial") ---------------------------------------------------
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "one-vehicle collision"
$ws.Range("D109").Value = "Yes"
$ws.Range("E109").Value = "Yes"
$ws.Range("F109").Value = "Yes"
$ws.Range("G109").Value = "Yes"
$ws.Range("K109").Value = "Yes"

# --- Row 111 ("55trial") ---------------------------------------------------
$ws.Range("B111").Value = "Unknown"
$ws.Range("C111").Value = "Unknown"
$ws.Range("D111").Value = "Unknown"
$ws.Range("E111").Value = "Unknown"
$ws.Range("F111").Value = "Unknown"
$ws.Range("G111").Value = "Unknown"
$ws.Range("H111").Value = "sedan, suv"

# --- Row 113 ("56trial") ---------------------------------------------------
$ws.Range("B113").Value = "N/A"
$ws.Range("C113").Value = "N/A"
$ws.Range("D113").Value = "No"
$ws.Range("E113").Value = "No"
$ws.Range("F113").Value = "No"
$ws.Range("G113").Value = "No"
$ws.Range("H113").Value = "N/A"
$ws.Range("I113").Value = "No"
$ws.Range("L113").Value = "No accident occurred"

# --- Row 117 ("58trial") ---------------------------------------------------
$ws.Range("B117").Value = 1
$ws.Range("C117").Value = "Unknown"
$ws.Range("D117").Value = "Yes"
$ws.Range("E117").Value = "Unknown"
$ws.Range("F117").Value = "No"
$ws.Range("G117").Value = "Yes"
$ws.Range("H117").Value = "Motorcycle or scooter?"
$ws.Range("I117").Value = "No"
$ws.Range("J117").Value = "Night"
$ws.Range("K117").Value = "Yes"

# --- Row 121 ("60trial") ---------------------------------------------------
$ws.Range("B121").Value = 2
$ws.Range("C121").Value = "T-bone"
$ws.Range("D121").Value = "Possibly"
$ws.Range("E121").Value = "Yes"
$ws.Range("F121").Value = "No"
$ws.Range("G121").Value = "Yes"
$ws.Range("I121").Value = "No"

# --- Row 127 ("63trial") ---------------------------------------------------
$ws.Range("B127").Value = 2
$ws.Range("C127").Value = "rear-end"
$ws.Range("D127").Value = "Yes"
$ws.Range("E127").Value = "Yes"
$ws.Range("F127").Value = "No"
$ws.Range("G127").Value = "Yes"

# --- Row 129 ("64trial") ---------------------------------------------------
$ws.Range("B129").Value = 1
$ws.Range("C129").Value = "sudden overturn"
$ws.Range("D129").Value = "Yes"
$ws.Range("E129").Value = "Yes"
$ws.Range("F129").Value = "Yes"
$ws.Range("G129").Value = "Yes"

# --- Row 135 ("67trial") ---------------------------------------------------
$ws.Range("B135").Value = "N/A"
$ws.Range("C135").Value = "N/A"
$ws.Range("D135").Value = "No"
$ws.Range("E135").Value = "No"
$ws.Range("F135").Value = "No"
$ws.Range("G135").Value = "No"
$ws.Range("H135").Value = "N/A"
$ws.Range("I135").Value = "No"
$ws.Range("K135").Value = "No"
$ws.Range("L135").Value = "Urban, pedrestrian crossing, no accident occurred"

# --- Row 139 ("69trial") ---------------------------------------------------
$ws.Range("B139").Value = 2
$ws.Range("C139").Value = "t-bone (train)"
$ws.Range("D139").Value = "Yes"
$ws.Range("E139").Value = "Yes"
$ws.Range("F139").Value = "Yes"
$ws.Range("G139").Value = "Yes"
$ws.Range("H139").Value = "truck, train"
$ws.Range("I139").Value = "No"
$ws.Range("J139").Value = "Day, clear"
$ws.Range("K139").Value = "No"
$ws.Range("L139").Value = "AI gets this one spot-on"

# --- Row 141 ("70trial") ---------------------------------------------------
$ws.Range("B141").Value = 1
$ws.Range("C141").Value = "Loss of control"
$ws.Range("D141").Value = "Unknown"
$ws.Range("E141").Value = "Unknown"
$ws.Range("F141").Value = "No"
$ws.Range("G141").Value = "Yes"
$ws.Range("H141").Value = "suv"
$ws.Range("I141").Value = "No"
$ws.Range("J141").Value = "Day, clear"
$ws.Range("K141").Value = "No"
$ws.Range("L141").Value = "Parking lot"

# --- Row 143 ("71trial") ---------------------------------------------------
$ws.Range("B143").Value = 2
$ws.Range("C143").Value = "t-bone"
$ws.Range("D143").Value = "Yes"
$ws.Range("E143").Value = "Yes"
$ws.Range("F143").Value = "Yes"
$ws.Range("G143").Value = "Yes"
$ws.Range("H143").Value = "suv, sedan"
$ws.Range("I143").Value = "No"
$ws.Range("J143").Value = "Day, clear"
$ws.Range("K143").Value = "No"

# --- Row 145 ("72trial") ---------------------------------------------------
$ws.Range("B145").Value = 2
$ws.Range("C145").Value = "t-bone"
$ws.Range("D145").Value = "Yes"
$ws.Range("E145").Value = "No"
$ws.Range("F145").Value = "Yes"
$ws.Range("G145").Value = "Yes"
$ws.Range("H145").Value = "suv, sedan"
$ws.Range("I145").Value = "No"
$ws.Range("J145").Value = "Day, clear"
$ws.Range("K145").Value = "No"

# --- Row 149 ("74trial") ---------------------------------------------------
$ws.Range("B149").Value = "2 or 3"
$ws.Range("C149").Value = "Unknown"
$ws.Range("D149").Value = "Unknown"
$ws.Range("E149").Value = "Unknown"
$ws.Range("F149").Value = "No"
$ws.Range("G149").Value = "Yes"
$ws.Range("K149").Value = "Yes"

# --- Row 153 ("76trial") ---------------------------------------------------
$ws.Range("B153").Value = 2
$ws.Range("C153").Value = "t-bone"
$ws.Range("D153").Value = "Yes"
$ws.Range("E153").Value = "Yes"
$ws.Range("F153").Value = "No"
$ws.Range("G153").Value = "Yes"
$ws.Range("H153").Value = "sedan, suv"
$ws.Range("I153").Value = "No"
$ws.Range("J153").Value = "Day, clear"
$ws.Range("K153").Value = "No"

# --- Small corrections to already-filled rows -------------------------------
# Row 120 ("60"): capitalize "yes" -> "Yes"
$ws.Range("E120").Value = "Yes"
# Row 124 ("62"): capitalize "yes" -> "Yes"
$ws.Range("E124").Value = "Yes"
# Row 138 ("69"): correct ambulance/firetruck/police need flags No -> Yes
$ws.Range("D138").Value = "Yes"
$ws.Range("E138").Value = "Yes"
$ws.Range("F138").Value = "Yes"

# --- Update the active selection to match where the author ended up --------
$ws.Range("H35").Select()
